{"js": "// Insert the new intro copy (empty spacer paragraph, two body paragraphs,\n// and a \"6 Separating Files\" Heading 1) right after the \"Write Up\" title,\n// and before the pre-existing trailing empty paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the first paragraph after the title - i.e. the first of the\n// existing (empty) trailing paragraphs. All new content is inserted\n// immediately before it, so it ends up right after \"Write Up\".\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Write Up\") {\n    anchor = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find paragraph following \"Write Up\"');\n}\n\n// Blank spacer paragraph.\nanchor.insertParagraph(\"\", Word.InsertLocation.before);\n\n// Intro paragraph.\nanchor.insertParagraph(\n  \"This week, we will be taking a look at how we can separate our files into individual modules. This will make things work faster and better. In addition, it can be very helpful for large projects where you have several different developers building the app. This way each employee can take their own piece of the pie (so-to-speak) and piece it back together to make a better whole.\",\n  Word.InsertLocation.before\n);\n\n// Lead-in paragraph.\nanchor.insertParagraph(\n  \"So, if this sounds like something that would interest you, then please join us for our new article entitled:\",\n  Word.InsertLocation.before\n);\n\n// New article heading.\nconst headingPara = anchor.insertParagraph(\"6 Separating Files\", Word.InsertLocation.before);\nheadingPara.style = \"Heading 1\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$cr = [char]13\n\n# Locate the \"Write Up\" title paragraph dynamically.\n$titleIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd($cr) -eq \"Write Up\") {\n        $titleIndex = $i\n        break\n    }\n}\nif ($titleIndex -eq -1) {\n    throw \"Could not find 'Write Up' paragraph\"\n}\n\n# The new content is inserted right after the title, i.e. right before\n# the paragraph that currently follows it (the first pre-existing empty\n# paragraph). Collapsing a range to the start of that paragraph and using\n# InsertBefore with embedded paragraph marks (chr 13) creates clean new\n# paragraphs that do not inherit the Title style.\n$anchorPara = $d.Paragraphs.Item($titleIndex + 1)\n$rng = $anchorPara.Range\n$rng.Collapse(1)  # wdCollapseStart\n\n$introText = \"This week, we will be taking a look at how we can separate our files into individual modules. This will make things work faster and better. In addition, it can be very helpful for large projects where you have several different developers building the app. This way each employee can take their own piece of the pie (so-to-speak) and piece it back together to make a better whole.\"\n$leadInText = \"So, if this sounds like something that would interest you, then please join us for our new article entitled:\"\n$headingText = \"6 Separating Files\"\n\n$rng.InsertBefore($cr + $introText + $cr + $leadInText + $cr + $headingText + $cr)\n\n# The heading paragraph is now the 4th new paragraph after the title\n# (blank spacer, intro, lead-in, heading).\n$headingIndex = $titleIndex + 4\n$headingPara = $d.Paragraphs.Item($headingIndex)\n$headingPara.Range.Style = \"Heading 1\"\n"}
